# Scheduled-runner update: refresh market-board derived figures
# (currentAveragePrice / *NQ / *HQ / LevePriceNQ / LevePriceHQ /
#  LeveProfitNQ / LeveProfitHQ, columns H-N) across the ALC, ARM, BSM,
# CRP, CUL, GSM, LTW and WVR sheets to their latest sampled values.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1381.32
$ws.Range("I132").Value = 1072.9524
$ws.Range("K132").Value = 3218.857199999999
$ws.Range("M132").Value = -688.8571999999995

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3958.1765
$ws.Range("J138").Value = 4642.4287
$ws.Range("L138").Value = 13927.2861
$ws.Range("N138").Value = -24207.2861

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1021
$ws.Range("I2").Value = 882.3684
$ws.Range("K2").Value = 882.3684
$ws.Range("M2").Value = -769.3684

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 802018.75
$ws.Range("I32").Value = 930758.9
$ws.Range("J32").Value = 13485.75
$ws.Range("K32").Value = 930758.9
$ws.Range("L32").Value = 13485.75
$ws.Range("M32").Value = -930471.9
$ws.Range("N32").Value = -14059.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2355.0625
$ws.Range("I45").Value = 1906.2307
$ws.Range("J45").Value = 4300
$ws.Range("K45").Value = 1906.2307
$ws.Range("L45").Value = 4300
$ws.Range("M45").Value = -1529.2307
$ws.Range("N45").Value = -5054

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5007507
$ws.Range("I61").Value = 9394
$ws.Range("J61").Value = 20001848
$ws.Range("K61").Value = 9394
$ws.Range("L61").Value = 20001848
$ws.Range("M61").Value = -9182
$ws.Range("N61").Value = -20002272

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2328021
$ws.Range("I74").Value = 3095027.5
$ws.Range("K74").Value = 3095027.5
$ws.Range("M74").Value = -3094153.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2328021
$ws.Range("I77").Value = 3095027.5
$ws.Range("K77").Value = 15475137.5
$ws.Range("M77").Value = -15470769.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1021
$ws.Range("I116").Value = 882.3684
$ws.Range("K116").Value = 882.3684
$ws.Range("M116").Value = 1411.6316

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 5007507
$ws.Range("I136").Value = 9394
$ws.Range("J136").Value = 20001848
$ws.Range("K136").Value = 28182
$ws.Range("L136").Value = 60005544
$ws.Range("M136").Value = -25632
$ws.Range("N136").Value = -60010644

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1021
$ws.Range("I3").Value = 882.3684
$ws.Range("K3").Value = 882.3684
$ws.Range("M3").Value = -768.3684

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H33").Value = 12024
$ws.Range("J33").Value = 12024
$ws.Range("L33").Value = 12024
$ws.Range("N33").Value = -12696

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 15476.883
$ws.Range("I105").Value = 13801.714
$ws.Range("J105").Value = 16649.5
$ws.Range("K105").Value = 13801.714
$ws.Range("L105").Value = 16649.5
$ws.Range("M105").Value = -12054.714
$ws.Range("N105").Value = -20143.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 7578734
$ws.Range("I134").Value = 3375.2222
$ws.Range("J134").Value = 41667850
$ws.Range("K134").Value = 10125.6666
$ws.Range("L134").Value = 125003550
$ws.Range("M134").Value = -7590.6666
$ws.Range("N134").Value = -125008620

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3289.923
$ws.Range("I16").Value = 4219.1113
$ws.Range("J16").Value = 1199.25
$ws.Range("K16").Value = 4219.1113
$ws.Range("L16").Value = 1199.25
$ws.Range("M16").Value = -3932.1113
$ws.Range("N16").Value = -1773.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5850108.5
$ws.Range("I31").Value = 6581342
$ws.Range("J31").Value = 243
$ws.Range("K31").Value = 6581342
$ws.Range("L31").Value = 243
$ws.Range("M31").Value = -6581047
$ws.Range("N31").Value = -833

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 5850108.5
$ws.Range("I34").Value = 6581342
$ws.Range("J34").Value = 243
$ws.Range("K34").Value = 6581342
$ws.Range("L34").Value = 243
$ws.Range("M34").Value = -6581140
$ws.Range("N34").Value = -647

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1072.3684
$ws.Range("J107").Value = 1628
$ws.Range("L107").Value = 1628
$ws.Range("N107").Value = -5468

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 3289.923
$ws.Range("I113").Value = 4219.1113
$ws.Range("J113").Value = 1199.25
$ws.Range("K113").Value = 4219.1113
$ws.Range("L113").Value = 1199.25
$ws.Range("M113").Value = -2049.1113
$ws.Range("N113").Value = -5539.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 7184681
$ws.Range("I4").Value = 18231484
$ws.Range("J4").Value = 4259.25
$ws.Range("K4").Value = 54694452
$ws.Range("L4").Value = 12777.75
$ws.Range("M4").Value = -54694340
$ws.Range("N4").Value = -13001.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2546.9333
$ws.Range("J34").Value = 2755.3333
$ws.Range("L34").Value = 8265.999899999999
$ws.Range("N34").Value = -8433.999899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 3974.7568
$ws.Range("I131").Value = 2328.1667
$ws.Range("J131").Value = 4293.4517
$ws.Range("K131").Value = 6984.500100000001
$ws.Range("L131").Value = 12880.3551
$ws.Range("M131").Value = -1944.500100000001
$ws.Range("N131").Value = -22960.3551

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 8061.8823
$ws.Range("I137").Value = 3007.8572
$ws.Range("J137").Value = 11599.7
$ws.Range("K137").Value = 9023.571599999999
$ws.Range("L137").Value = 34799.10000000001
$ws.Range("M137").Value = -3923.571599999999
$ws.Range("N137").Value = -44999.10000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 17750.824
$ws.Range("I70").Value = 26199.7
$ws.Range("J70").Value = 5681
$ws.Range("K70").Value = 26199.7
$ws.Range("L70").Value = 5681
$ws.Range("M70").Value = -25929.7
$ws.Range("N70").Value = -6221

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 17750.824
$ws.Range("I73").Value = 26199.7
$ws.Range("J73").Value = 5681
$ws.Range("K73").Value = 26199.7
$ws.Range("L73").Value = 5681
$ws.Range("M73").Value = -25263.7
$ws.Range("N73").Value = -7553

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2749.5
$ws.Range("J80").Value = 4000
$ws.Range("L80").Value = 4000
$ws.Range("N80").Value = -5996

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2749.5
$ws.Range("J83").Value = 4000
$ws.Range("L83").Value = 20000
$ws.Range("N83").Value = -29984

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 13004.223
$ws.Range("I102").Value = 13067.25
$ws.Range("K102").Value = 13067.25
$ws.Range("M102").Value = -11445.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 13825.3125
$ws.Range("I132").Value = 11188.3
$ws.Range("J132").Value = 18220.334
$ws.Range("K132").Value = 33564.89999999999
$ws.Range("L132").Value = 54661.00199999999
$ws.Range("M132").Value = -31034.89999999999
$ws.Range("N132").Value = -59721.00199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3356.0908
$ws.Range("I46").Value = 950
$ws.Range("K46").Value = 950
$ws.Range("M46").Value = -762

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 4951.731
$ws.Range("I93").Value = 4099.0527
$ws.Range("J93").Value = 7266.143
$ws.Range("K93").Value = 4099.0527
$ws.Range("L93").Value = 7266.143
$ws.Range("M93").Value = -2851.0527
$ws.Range("N93").Value = -9762.143

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 5787.375
$ws.Range("I100").Value = 5875.5
$ws.Range("J100").Value = 5699.25
$ws.Range("K100").Value = 5875.5
$ws.Range("L100").Value = 5699.25
$ws.Range("M100").Value = -5334.5
$ws.Range("N100").Value = -6781.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4805.8887
$ws.Range("I122").Value = 3351.2
$ws.Range("K122").Value = 10053.6
$ws.Range("M122").Value = -7603.599999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 4000
$ws.Range("I96").Value = 4000
$ws.Range("K96").Value = 4000
$ws.Range("M96").Value = -2627

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1436.125
$ws.Range("I100").Value = 698
$ws.Range("K100").Value = 1396
$ws.Range("M100").Value = -855

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()
